$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix typo: "CEIP workign group..." -> "CEIP working group..."
$ws.Range("B12").Value = "CEIP working group on U.S.-Russia relations "

# Fix typo: "Dartmouth Confernece " -> "Dartmouth Conference "
$ws.Range("B3").Value = "Dartmouth Conference "

# Replace short "CISAC Russia Dialogue" with the full official name
$ws.Range("A4").Value = "National Academy of Sciences Committee on International Security and Arms Control Russia Dialogue"

# Update start date for the Ivanov-Talbott-Albright Dialogue row
$ws.Range("D37").Value = 2009

# Reset the view: scroll to top-left and select A4
$ws.Activate()
$ws.Range("A4").Select()
